$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 28410.104
$ws.Range("I64").Value = 102410
$ws.Range("J64").Value = 2892.8965
$ws.Range("K64").Value = 102410
$ws.Range("L64").Value = 2892.8965
$ws.Range("M64").Value = -102162
$ws.Range("N64").Value = -3388.8965
$ws.Range("H67").Value = 28410.104
$ws.Range("I67").Value = 102410
$ws.Range("J67").Value = 2892.8965
$ws.Range("K67").Value = 102410
$ws.Range("L67").Value = 2892.8965
$ws.Range("M67").Value = -101552
$ws.Range("N67").Value = -4608.8965
$ws.Range("H76").Value = 3324.5
$ws.Range("I76").Value = 3048.7144
$ws.Range("K76").Value = 3048.7144
$ws.Range("M76").Value = -2733.7144
$ws.Range("H79").Value = 3324.5
$ws.Range("I79").Value = 3048.7144
$ws.Range("K79").Value = 3048.7144
$ws.Range("M79").Value = -1956.7144
$ws.Range("H86").Value = 7418.6665
$ws.Range("I86").Value = 1746.6666
$ws.Range("J86").Value = 13090.667
$ws.Range("K86").Value = 1746.6666
$ws.Range("L86").Value = 13090.667
$ws.Range("M86").Value = -623.6666
$ws.Range("N86").Value = -15336.667
$ws.Range("H89").Value = 7418.6665
$ws.Range("I89").Value = 1746.6666
$ws.Range("J89").Value = 13090.667
$ws.Range("K89").Value = 8733.333000000001
$ws.Range("L89").Value = 65453.335
$ws.Range("M89").Value = -3117.333000000001
$ws.Range("N89").Value = -76685.33499999999
$ws.Range("H125").Value = 1949.5555
$ws.Range("I125").Value = 4344
$ws.Range("J125").Value = 1470.6666
$ws.Range("K125").Value = 39096
$ws.Range("L125").Value = 13235.9994
$ws.Range("M125").Value = -36636
$ws.Range("N125").Value = -18155.9994
$ws.Range("H129").Value = 878.94446
$ws.Range("I129").Value = 526.4286
$ws.Range("J129").Value = 964.0345
$ws.Range("K129").Value = 1579.2858
$ws.Range("L129").Value = 2892.1035
$ws.Range("M129").Value = 3420.7142
$ws.Range("N129").Value = -12892.1035
$ws.Range("H135").Value = 993.1429000000001
$ws.Range("J135").Value = 800
$ws.Range("L135").Value = 7200
$ws.Range("N135").Value = -12270

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100.25
$ws.Range("I4").Value = 67
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 67
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 49
$ws.Range("N4").Value = -432
$ws.Range("H6").Value = 7000.375
$ws.Range("I6").Value = 5001
$ws.Range("J6").Value = 8999.75
$ws.Range("K6").Value = 5001
$ws.Range("L6").Value = 8999.75
$ws.Range("M6").Value = -4828
$ws.Range("N6").Value = -9345.75
$ws.Range("H9").Value = 10249.75
$ws.Range("J9").Value = 10249.75
$ws.Range("L9").Value = 10249.75
$ws.Range("N9").Value = -10589.75
$ws.Range("H20").Value = 10249.75
$ws.Range("J20").Value = 10249.75
$ws.Range("L20").Value = 10249.75
$ws.Range("N20").Value = -10789.75
$ws.Range("H88").Value = 2460.4614
$ws.Range("I88").Value = 2568.6
$ws.Range("J88").Value = 2100
$ws.Range("K88").Value = 2568.6
$ws.Range("L88").Value = 2100
$ws.Range("M88").Value = -2162.6
$ws.Range("N88").Value = -2912
$ws.Range("H91").Value = 2460.4614
$ws.Range("I91").Value = 2568.6
$ws.Range("J91").Value = 2100
$ws.Range("K91").Value = 2568.6
$ws.Range("L91").Value = 2100
$ws.Range("M91").Value = -1164.6
$ws.Range("N91").Value = -4908

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 124368.11
$ws.Range("I86").Value = 221301.2
$ws.Range("K86").Value = 221301.2
$ws.Range("M86").Value = -220178.2
$ws.Range("H89").Value = 124368.11
$ws.Range("I89").Value = 221301.2
$ws.Range("K89").Value = 1106506
$ws.Range("M89").Value = -1100890

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6542.6665
$ws.Range("I5").Value = 983.5
$ws.Range("K5").Value = 2950.5
$ws.Range("M5").Value = -2838.5
$ws.Range("H122").Value = 21025.6
$ws.Range("J122").Value = 34599.332
$ws.Range("L122").Value = 311393.988
$ws.Range("N122").Value = -316293.988
$ws.Range("H131").Value = 828.63635
$ws.Range("I131").Value = 614
$ws.Range("J131").Value = 844.9674
$ws.Range("K131").Value = 1842
$ws.Range("L131").Value = 2534.9022
$ws.Range("M131").Value = 3198
$ws.Range("N131").Value = -12614.9022
$ws.Range("H135").Value = 6542.6665
$ws.Range("I135").Value = 983.5
$ws.Range("K135").Value = 8851.5
$ws.Range("M135").Value = -6316.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 95718.27
$ws.Range("I70").Value = 204410.8
$ws.Range("J70").Value = 5141.1665
$ws.Range("K70").Value = 204410.8
$ws.Range("L70").Value = 5141.1665
$ws.Range("M70").Value = -204140.8
$ws.Range("N70").Value = -5681.1665
$ws.Range("H73").Value = 95718.27
$ws.Range("I73").Value = 204410.8
$ws.Range("J73").Value = 5141.1665
$ws.Range("K73").Value = 204410.8
$ws.Range("L73").Value = 5141.1665
$ws.Range("M73").Value = -203474.8
$ws.Range("N73").Value = -7013.1665
$ws.Range("H80").Value = 3832.5
$ws.Range("I80").Value = 3740
$ws.Range("J80").Value = 3925
$ws.Range("K80").Value = 3740
$ws.Range("L80").Value = 3925
$ws.Range("M80").Value = -2742
$ws.Range("N80").Value = -5921
$ws.Range("H83").Value = 3832.5
$ws.Range("I83").Value = 3740
$ws.Range("J83").Value = 3925
$ws.Range("K83").Value = 18700
$ws.Range("L83").Value = 19625
$ws.Range("M83").Value = -13708
$ws.Range("N83").Value = -29609

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1584.5714
$ws.Range("J7").Value = 2420
$ws.Range("L7").Value = 2420
$ws.Range("N7").Value = -2644
$ws.Range("H16").Value = 78323.84
$ws.Range("I16").Value = 100620.9
$ws.Range("K16").Value = 100620.9
$ws.Range("M16").Value = -100450.9
$ws.Range("H61").Value = 2492.1428
$ws.Range("I61").Value = 1862.0625
$ws.Range("J61").Value = 4508.4
$ws.Range("K61").Value = 1862.0625
$ws.Range("L61").Value = 4508.4
$ws.Range("M61").Value = -1660.0625
$ws.Range("N61").Value = -4912.4
$ws.Range("H113").Value = 2492.1428
$ws.Range("I113").Value = 1862.0625
$ws.Range("J113").Value = 4508.4
$ws.Range("K113").Value = 1862.0625
$ws.Range("L113").Value = 4508.4
$ws.Range("M113").Value = 307.9375
$ws.Range("N113").Value = -8848.4
$ws.Range("H126").Value = 1584.5714
$ws.Range("J126").Value = 2420
$ws.Range("L126").Value = 7260
$ws.Range("N126").Value = -12200

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1997.3077
$ws.Range("I136").Value = 677.2727
$ws.Range("J136").Value = 2965.3333
$ws.Range("K136").Value = 2031.8181
$ws.Range("L136").Value = 8895.999899999999
$ws.Range("M136").Value = 518.1819
$ws.Range("N136").Value = -13995.9999
